$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two data rows that were dropped from the sheet:
#  - old row 26 "RM 232"
#  - old row 28 "SC 92" (becomes row 27 once row 26 is removed)
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Update the "F" column values that changed for the surviving rows.
$ws.Range("F19").Value = 17.81
$ws.Range("F21").Value = ""
$ws.Range("F23").Value = 16.48
$ws.Range("F27").Value = ""
$ws.Range("F33").Value = 17.53
